# Remove the 2005-2009 data rows (original rows 2-6) from Sheet1.
# The remaining 2010-2013 rows shift up to rows 2-5, matching the target
# layout: dimension A1:J10 -> A1:J5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:6").Delete()
